# Todo.xlsx: add a new PBI/Task row for "[ENGINE] Object::Hide() function"
# just above the "Add Imaterial / randomly rotated poisson disk" block
# (i.e. directly under the existing "Allow creation of primitive objects:
# {Box, Sphere}" row), and move the selection to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 14, pushing the existing row 14 ("Allow creation
# of primitive objects: {Box, Sphere}") and everything below it down by one.
$ws.Rows("14:14").Insert()

# Populate the freshly inserted row with the new task text.
$ws.Range("B14").Value = "[ENGINE] Object::Hide() function"

# Match the author's final cursor position recorded in the sheet view.
$ws.Range("D11").Select()
